# New crime data collected - weekly CompStat report update (22nd Precinct / Central Park)
#
# Updates:
#  - Report header: volume/week number and the "week covering" date range
#  - Weekly "Murder" (row16) and "Rape" (row17) rows: week-to-date counts went to 0,
#    which this report renders as the text placeholders "0"/"***.*" rather than numbers
#  - "Gr. Larceny" (row19): week-to-date counts went from 0 (text placeholders) to real
#    numbers, plus downstream 28-day/YTD/trend figures were recalculated
#  - TOTAL row (row21) figures recalculated
#  - "Petit Larceny" (row24), "Misd. Assault" (row26) and "Hate Crimes" (row31):
#    same "text placeholder" <-> "real number" swaps plus recalculated figures
#
# Note: for cells that swap between the text placeholder ("0" / "***.*") representation
# and a real numeric representation, we use Range.Copy() from a same-shaped, untouched
# donor cell elsewhere on the sheet so that both the cell's value/type AND its number
# format (style) flip correctly together, exactly as Excel does when you retype a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header (A8: volume/number banner, C9: week-covering date range) ---
$ws.Range("A8").Value = "Volume 32   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/8/2025  Through  12/14/2025"

# --- Row 16 (Murder): week-to-date counts now 0 -> shown as text placeholders ---
$ws.Range("C14").Copy($ws.Range("C16"))     # 1 (number)   -> "0" (text placeholder)
$ws.Range("C14").Copy($ws.Range("D16"))     # 1 (number)   -> "0" (text placeholder)
$ws.Range("E14").Copy($ws.Range("E16"))     # 0 (number)   -> "***.*" (text placeholder)
$ws.Range("N16").Value = -94.38775510204    # 32-year trend recalculated

# --- Row 17 (Rape): week-to-date counts now 0 -> shown as text placeholders ---
$ws.Range("C14").Copy($ws.Range("C17"))     # 1 (number)   -> "0" (text placeholder)
$ws.Range("C14").Copy($ws.Range("G17"))     # 1 (number)   -> "0" (text placeholder)
$ws.Range("E14").Copy($ws.Range("H17"))     # 0 (number)   -> "***.*" (text placeholder)
$ws.Range("L17").Value = 122.222222222222   # 2-year trend recalculated

# --- Row 19 (Gr. Larceny): week-to-date placeholders now real counts ---
$ws.Range("J18").Copy($ws.Range("C19"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("J18").Copy($ws.Range("D19"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("K15").Copy($ws.Range("E19"))     # "***.*" (text placeholder) -> 0 (number)
$ws.Range("J18").Copy($ws.Range("F19"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("G19").Value = 3                  # 28-day 2024 count
$ws.Range("H19").Value = -66.666666666666   # 28-day % change
$ws.Range("I19").Value = 44                 # YTD 2025 count
$ws.Range("J19").Value = 44                 # YTD 2024 count
$ws.Range("L19").Value = -8.333333333333    # 2-year trend
$ws.Range("M19").Value = -29.032258064516   # 15-year trend
$ws.Range("N19").Value = -74.71264367816    # 32-year trend

# --- Row 21 (TOTAL): recalculated after the above changes ---
$ws.Range("C21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 80
$ws.Range("J21").Value = 99
$ws.Range("K21").Value = -19.191919191919
$ws.Range("L21").Value = -4.761904761904
$ws.Range("M21").Value = -19.191919191919
$ws.Range("N21").Value = -82.721382289416

# --- Row 24 (Petit Larceny): week-to-date placeholders now real counts ---
$ws.Range("J18").Copy($ws.Range("D24"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("N14").Copy($ws.Range("E24"))     # "***.*" (text placeholder) -> -100 (number)
$ws.Range("J24").Value = 38                 # YTD 2024 count
$ws.Range("K24").Value = -2.631578947368    # YTD % change

# --- Row 26 (Misd. Assault): week-to-date counts now 0 -> text placeholders ---
$ws.Range("C14").Copy($ws.Range("D26"))     # 1 (number) -> "0" (text placeholder)
$ws.Range("E14").Copy($ws.Range("E26"))     # -100 (number) -> "***.*" (text placeholder)

# --- Row 31 (Hate Crimes): week-to-date placeholders now real counts ---
$ws.Range("J18").Copy($ws.Range("D31"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("N14").Copy($ws.Range("E31"))     # "***.*" (text placeholder) -> -100 (number)
$ws.Range("J18").Copy($ws.Range("G31"))     # "0" (text placeholder) -> 1 (number)
$ws.Range("N14").Copy($ws.Range("H31"))     # "***.*" (text placeholder) -> -100 (number)
$ws.Range("J31").Value = 4                  # YTD 2024 count
$ws.Range("K31").Value = -50                # YTD % change
